$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "已完成" note next to the mx/全程跑实验 row (row 4), in column E.
$ws.Range("E4").Value = "已完成"

# Move/park the selection where the author left it before saving.
$ws.Range("F5").Select() | Out-Null
